$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.68"
$ws.Range("E2").Value = "'6.42%"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'27.28"
$ws.Range("E3").Value = "'-0.91%"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'4.797"
$ws.Range("E4").Value = "'2.09%"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.06252"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'6.837"
$ws.Range("E6").Value = "'1.52%"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'0.8785"
$ws.Range("E7").Value = "'3.30%"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'0.9430"
$ws.Range("E8").Value = "'3.55%"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.1461"
$ws.Range("E9").Value = "'3.84%"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.05011"
$ws.Range("E10").Value = "'4.08%"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.07270"
$ws.Range("E11").Value = "'2.70%"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.03154"
$ws.Range("E12").Value = "'0.94%"
$ws.Range("G12").Value = "'3"
$ws.Range("D13").Value = "'0.09036"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("G13").Value = "'3"
$ws.Range("D14").Value = "'0.001546"
$ws.Range("E14").Value = "'0.92%"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'0.0006282"
$ws.Range("E15").Value = "'1.94%"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.005746"
$ws.Range("E16").Value = "'-4.58%"
$ws.Range("G16").Value = "'3"
$ws.Range("D17").Value = "'3.471"
$ws.Range("E17").Value = "'0.55%"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'3.267"
$ws.Range("E18").Value = "'3.19%"
$ws.Range("G18").Value = "'3"
$ws.Range("D19").Value = "'2.235"
$ws.Range("E19").Value = "'3.20%"
$ws.Range("G19").Value = "'3"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.1310"
$ws.Range("E21").Value = "'-0.02%"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'3.856"
$ws.Range("E22").Value = "'-5.77%"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'0.04311"
$ws.Range("E23").Value = "'1.40%"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'0.001176"
$ws.Range("E24").Value = "'-3.60%"
$ws.Range("G24").Value = "'3"
$ws.Range("D25").Value = "'0.004262"
$ws.Range("E25").Value = "'4.47%"
$ws.Range("G25").Value = "'3"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.14%"
$ws.Range("G26").Value = "'3"
$ws.Range("D27").Value = "'0.0001614"
$ws.Range("E27").Value = "'-1.60%"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("G38").Value = "'3"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.04020"
$ws.Range("E40").Value = "'3.82%"
$ws.Range("G40").Value = "'3"
$ws.Range("D41").Value = "'0.006366"
$ws.Range("E41").Value = "'54.80%"
$ws.Range("G41").Value = "'3"
$ws.Range("D42").Value = "'0.1149"
$ws.Range("E42").Value = "'3.46%"
$ws.Range("G42").Value = "'3"
$ws.Range("D43").Value = "'0.002209"
$ws.Range("E43").Value = "'2.86%"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.01355"
$ws.Range("E44").Value = "'3.17%"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.00005132"
$ws.Range("E45").Value = "'0.00%"
$ws.Range("G45").Value = "'3"
$ws.Range("E46").Value = "'-0.16%"
$ws.Range("G46").Value = "'3"
$ws.Range("D47").Value = "'2.055"
$ws.Range("E47").Value = "'2,852.12%"
$ws.Range("G47").Value = "'3"
$ws.Range("E48").Value = "'-12.22%"
$ws.Range("G48").Value = "'3"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.16%"
$ws.Range("G49").Value = "'3"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("G50").Value = "'3"
$ws.Range("G51").Value = "'3"
